$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Retained Earnings, Total Equity, and Total Liabilities & Equity
# rows (11-13) for columns E, F, G with the corrected figures. Values are
# kept as text (matching the existing text-formatted figures in the sheet,
# e.g. "1,305") rather than being auto-converted to numbers.
$updates = @{
    "E11" = "3,912"
    "F11" = "5,558"
    "G11" = "7,543"
    "E12" = "4,012"
    "F12" = "5,658"
    "G12" = "7,643"
    "E13" = "4,251"
    "F13" = "5,908"
    "G13" = "7,903"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}

# Remove column H (the "units" header and its empty cells below) entirely,
# shrinking the sheet's used range/dimension from A1:H13 to A1:G13.
$ws.Columns("H").Delete()
